# Add a new "2022-Q1" sheet (fund holdings detail, same layout as the other
# quarterly sheets) positioned immediately before "总计", and update the
# "总计" (totals) sheet with a new leading row summarizing 2022-Q1.

$wb = $excel.ActiveWorkbook

# --- 1. Update the "总计" sheet FIRST ---------------------------------------
# NOTE: do this before inserting any new sheet. A sheet reference fetched via
# Worksheets.Item(name) tracks the sheet's *position*, not a stable identity
# - once a new sheet is spliced in before "总计" its tab position shifts, so
# a reference obtained beforehand would silently start pointing at the wrong
# sheet. Finishing all "总计" edits while it is still the last sheet avoids
# that entirely.
#
# Rewrite the data rows directly (rather than inserting a row) so no
# floating-point drift is introduced in the untouched numeric cells; row 2
# becomes the new 2022-Q1 summary and every following row shifts down by one.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("B2:B7").NumberFormat = "@"

$totalRows = @(
  @(0, "2022-Q1", 8,  0.52),
  @(1, "2021-Q4", 29, 10.39),
  @(2, "2021-Q3", 23, 7.82),
  @(3, "2021-Q2", 26, 18.41),
  @(4, "2021-Q1", 34, 16.04),
  @(5, "2020-Q4", 34, 11.76)
)

$r = 2
foreach ($row in $totalRows) {
    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Row 7 is brand-new (the sheet previously only went to row 6), so its
# column-A cell needs the same bold/bordered/centered styling the other
# column-A cells already carry. Copy the format from the cell above (A6)
# rather than re-declaring Font/Borders/Alignment by hand, since that hits
# a *different* (but visually similar) pre-existing style than the one the
# other column-A cells use.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)

# --- 2. Create the "2022-Q1" sheet ------------------------------------------
# Duplicate an existing quarterly sheet ("2021-Q4") so the new sheet inherits
# the exact same column layout/styles (bold header row, bold+bordered column
# A, text-typed B:G columns), then drop it right before "总计".
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$srcSheet.Copy($totalSheet)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The source sheet has 30 data rows; the new sheet only needs 9 (header + 8
# funds), so drop the extra rows.
$newSheet.Rows("10:30").Delete()

# Columns B:G hold text-like values (fund codes, names, percentages) that
# must stay text (not get coerced to numbers, which would e.g. drop the
# leading zero on fund codes like "005984").
$newSheet.Range("B2:G9").NumberFormat = "@"

$fundRows = @(
  @(0, "005984", "兴业聚华混合A",               "21.22", "24.30", "0.99", "0.2101", 9),
  @(1, "000587", "大成灵活配置混合",             "2.39",  "85.02", "3.87", "0.0925", 5),
  @(2, "005985", "兴业聚华混合C",               "7.68",  "24.30", "0.99", "0.0760", 9),
  @(3, "003601", "申万菱信安鑫精选混合A",         "4.50",  "24.28", "1.21", "0.0544", 6),
  @(4, "515860", "嘉实中证新兴科技100策略ETF",    "2.25",  "98.94", "2.36", "0.0531", 9),
  @(5, "004351", "汇丰晋信珠三角区域发展混合",     "0.51",  "93.92", "4.92", "0.0251", 5),
  @(6, "004730", "建信量化事件驱动股票",          "0.72",  "91.24", "1.36", "0.0098", 8),
  @(7, "003602", "申万菱信安鑫精选混合C",         "0.01",  "24.28", "1.21", "0.0001", 6)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

Write-Host "2022-Q1 sheet added and 总计 sheet updated"
